# NewEngineFormulas.xlsx — add 4 new "d13c" engine-level rows into the
# first block (First LEVEL / d16g... block) on Arkusz1, mirroring the
# existing d13c rows already present in the second block (Last Level).
#
# Net effect on the sheet: rows that used to live at 151-164 are pushed
# down by 4 (to 155-168), and four brand-new rows are created at 151-154
# holding the d13c420 / d13c460 / d13c500 / d13c540 data (same truck
# names/values as appear later at rows 160/161/163/164, but wired into
# the $E$130-anchored formula group that rows 131-149 belong to).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert 4 blank rows right above the current row 151. Everything at
#    151 and below (old 151-164) shifts down to 155-168; Excel updates
#    all relative formulas / MIN/MAX ranges / dimension automatically.
$ws.Rows("151:154").Insert()

# 2) Seed the formatting of the new rows from the (now-shifted) sibling
#    data rows directly below them, so styles/number formats match the
#    rest of the d16g/d13c table instead of inheriting the header row's
#    format that a raw row-insert leaves behind.
$ws.Range("A155:H158").Copy()
$ws.Range("A151:H154").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Truck name (col A), rated torque input (col B) and the "Acc EHP"
#    column G for the four new rows.
$ws.Range("A151").Value = "d13c420.sii"
$ws.Range("B151").Value = 420
$ws.Range("G151").Value = 0

$ws.Range("A152").Value = "d13c460.sii"
$ws.Range("B152").Value = 460
$ws.Range("G152").Value = 2

$ws.Range("A153").Value = "d13c500.sii"
$ws.Range("B153").Value = 500
$ws.Range("G153").Value = 6

$ws.Range("A154").Value = "d13c540.sii"
$ws.Range("B154").Value = 540
$ws.Range("G154").Value = 11

# 4) Formulas for F151:F154 and H151:H154 — same shape as the other rows
#    in this block (131-149), which this new block extends: F anchors on
#    $E$130, H rounds the level from F using $C$4/$D$4/$C$2/$D$2.
for ($i = 151; $i -le 154; $i++) {
    $ws.Range("F$i").Formula = "=(`$E`$130+B$i)*(IF(B$i<=`$C`$8,1,0)*`$D`$8+IF(AND(B$i>`$C`$8,B$i<=`$C`$9),1,0)*`$D`$9+IF(AND(B$i>`$C`$9,B$i<=`$C`$10),1,0)*`$D`$10+IF(AND(B$i>`$C`$10,B$i<=`$C`$11),1,0)*`$D`$11+IF(AND(B$i>`$C`$11,B$i<=`$C`$12),1,0)*`$D`$12+IF(AND(B$i>`$C`$12,B$i<=`$C`$13),1,0)*`$D`$13+IF(AND(B$i>`$C`$13,B$i<=`$C`$14),1,0)*`$D`$14+IF(AND(B$i>`$C`$14,B$i<=`$C`$15),1,0)*`$D`$15+IF(AND(B$i>`$C`$15,B$i<=`$C`$16),1,0)*`$D`$16+IF(AND(B$i>`$C`$16,B$i<=`$C`$17),1,0)*`$D`$17)"
    $ws.Range("H$i").Formula = "=ROUND(((F$i-`$C`$4)/(`$D`$4-`$C`$4))*(`$D`$2-`$C`$2),0)+`$C`$2"
}

# 5) Restore the view: scrolled down toward the newly-edited area with
#    K154 selected.
$ws.Application.ActiveWindow.ScrollRow = 136
$ws.Range("K154").Select()

$wb.Application.Calculate()
